# Apply edits described by the commit:
# "Able to find the operator at correct location"
#
# 1. Update the C-column values for rows 2-10 from 193 to 188.
# 2. Move the active selection to N7 (so the correct location is
#    highlighted/selected when the workbook is reopened).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2:C10 values (193 -> 188)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 188
}

# Update the saved selection to N7
$ws.Range("N7").Select()
